$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"

# --- Row 16 ---
$ws.Range("C16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "0"
$ws.Range("C16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = "***.*"
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -55.555555555555

# --- Row 17 ---
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -71.428571428571
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = -55.555555555555
$ws.Range("L17").Value = -20

# --- Row 18 ---
$ws.Range("C18").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "0"
$ws.Range("C18").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = "***.*"
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -40
$ws.Range("F18").Copy()
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("I18").Value = 1
$ws.Range("K18").Value = -66.666666666666
$ws.Range("H18").Copy()
$ws.Range("L18").PasteSpecial(-4122)
$ws.Range("L18").Value = -50

# --- Row 19 ---
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -55.555555555555
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -9.375
$ws.Range("I19").Value = 15
$ws.Range("J19").Value = 15
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = -40

# --- Row 20 ---
$ws.Range("D20").Value = 2
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -20
$ws.Range("J20").Value = 3

# --- Row 21 ---
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -66.666666666666
$ws.Range("F21").Value = 56
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = -22.222222222222
$ws.Range("I21").Value = 20
$ws.Range("J21").Value = 32
$ws.Range("K21").Value = -37.5
$ws.Range("L21").Value = -51.219512195122

# --- Row 23 ---
$ws.Range("G23").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 1
$ws.Range("H23").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 2
$ws.Range("G23").Copy()
$ws.Range("J23").PasteSpecial(-4122)
$ws.Range("J23").Value = 1
$ws.Range("H23").Copy()
$ws.Range("K23").PasteSpecial(-4122)
$ws.Range("K23").Value = -100
$ws.Range("H23").Copy()
$ws.Range("L23").PasteSpecial(-4122)
$ws.Range("L23").Value = -100

# --- Row 24 ---
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -25.714285714285
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = 136
$ws.Range("H24").Value = -10.294117647058
$ws.Range("I24").Value = 61
$ws.Range("J24").Value = 59
$ws.Range("K24").Value = 3.389830508474
$ws.Range("L24").Value = 52.5

# --- Row 25 ---
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = -4.878048780487
$ws.Range("I25").Value = 18
$ws.Range("J25").Value = 19
$ws.Range("K25").Value = -5.263157894736
$ws.Range("L25").Value = 20

# --- Row 27 ---
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -85.714285714285
$ws.Range("D27").Copy()
$ws.Range("I27").PasteSpecial(-4122)
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = -75
$ws.Range("L27").Value = -80

# --- Row 30 ---
$ws.Range("F30").Copy()
$ws.Range("I30").PasteSpecial(-4122)
$ws.Range("I30").Value = 1
